# Updated test cases sheet & deleted old test recording
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: new test case TC2 for the "Your Account Page" section. Fill the new
# test title / type / section / id first (this is the order the new shared
# strings get interned in), then go back and fill the Priority column.
$ws.Range("C36").Value = "Verify all the options from your accounts page"
$ws.Range("D36").Value = "Medium "
$ws.Range("D36").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("E36").Value = "Sanity "
$ws.Range("A36").Value = "Your Account Page "
$ws.Range("B36").Value = "TC2"

# Row 35: fill in the previously-blank Priority (D) and Type (E) for the
# "Your Account Page" TC1 row. D35 previously had the default bordered style;
# Excel drops the border when the value is typed in here (matches the
# "Medium "/"High " cells used elsewhere that lack borders), so clear the
# border after setting the value.
$ws.Range("E35").Value = "Sanity "
$ws.Range("D35").Value = "High "
$ws.Range("D35").Borders.LineStyle = -4142   # xlLineStyleNone

# Move the active selection to reflect where the author left off editing.
$ws.Activate()
$ws.Range("C29").Select()
